$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("G5").Value = 2.45
$ws.Range("I5").Value = 3.1
$ws.Range("J5").Value = 3.25
$ws.Range("W5").Value = 6.5
$ws.Range("Z5").Value = 23
$ws.Range("AG5").Value = 1250
$ws.Range("AH5").Value = 8
$ws.Range("AI5").Value = 15
$ws.Range("AJ5").Value = 13
$ws.Range("AZ5").Value = 67

# Row 7
$ws.Range("Q7").Value = 1.75
$ws.Range("R7").Value = 2.05

# Row 14
$ws.Range("G14").Value = 3.75
$ws.Range("J14").Value = 3.75
$ws.Range("K14").Value = 2.75
$ws.Range("N14").Value = 29
$ws.Range("Q14").Value = 1.3
$ws.Range("R14").Value = 3.5
$ws.Range("S14").Value = 1.18
$ws.Range("T14").Value = 4.5
$ws.Range("W14").Value = 23
$ws.Range("X14").Value = 29
$ws.Range("Y14").Value = 15
$ws.Range("AF14").Value = 26
$ws.Range("AM14").Value = 15
$ws.Range("AT14").Value = 4.5
$ws.Range("AX14").Value = 9

# Row 22
$ws.Range("Q22").Value = 1.33

# Row 23
$ws.Range("R23").Value = 1.58

# Row 30
$ws.Range("M30").Value = 1.03
$ws.Range("O30").Value = 1.25

# Row 31
$ws.Range("M31").Value = 1.02
$ws.Range("O31").Value = 1.11

# Row 32
$ws.Range("M32").Value = 1.05
$ws.Range("O32").Value = 1.33

# Row 33
$ws.Range("M33").Value = 1.02
$ws.Range("O33").Value = 1.19

# Row 35
$ws.Range("J35").Value = 2.88
